$wb = $excel.ActiveWorkbook

# Updated "想去人数" (people interested) counts for a handful of events.
# Apply to both sheets that carry this data set: "展览" and "全部类型".
$updates = @{
    2  = 2925
    4  = 102
    5  = 6711
    6  = 1658
    7  = 19
    8  = 29
    9  = 55
    10 = 111
    11 = 23
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
